$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status updates: "Strip out Demo Application" and "Character Movement" moved
#     from NOT STARTED to IN PROGRESS ---
$ws.Range("F29").Value = "IN PROGRESS"
$ws.Range("F30").Value = "IN PROGRESS"

# --- Two new backlog items appended as rows 35 & 36 ---
$ws.Range("A35").Value = 30
$ws.Range("B35").Value = "Game Treatment 2"
$ws.Range("C35").Value = "Comic Book Style Storyboard "
$ws.Range("D35").Value = "-"
$ws.Range("E35").Value = "Andy"
$ws.Range("F35").Value = "NOT STARTED"

$ws.Range("A36").Value = 31
$ws.Range("C36").Value = "Provice Concept Art Backlog"
$ws.Range("B36").Value = "Art 8"
$ws.Range("D36").Value = "-"
$ws.Range("E36").Value = "Hiren"
$ws.Range("F36").Value = "NOT STARTED"

# --- Widen the title merge from A1:E1 to A1:F1 ---
$ws.Range("A1:E1").UnMerge()
$ws.Range("A1:F1").Merge()

# --- Extend the "Status" conditional formatting (NOT STARTED/IN PROGRESS/DONE)
#     so it keeps covering the new rows ---
$statusRange = $ws.Range("F4:F34")
for ($i = 1; $i -le $statusRange.FormatConditions.Count; $i++) {
    $statusRange.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("F4:F36"))
}

# --- Extend the data-validation dropdown list to cover the new rows ---
$ws.Range("F4:F34").Validation.Delete()
$ws.Range("F4:F36").Validation.Add(3, 1, 1, "=`$I`$3:`$I`$5")

# --- Restore the cursor position as recorded in the saved workbook ---
$null = $ws.Range("C16").Select()
